$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 1171
$ws1.Range("F3").Value = 1084
$ws1.Range("F4").Value = 1887
$ws1.Range("F5").Value = 592
$ws1.Range("F6").Value = 1235
$ws1.Range("F10").Value = 325
$ws1.Range("F11").Value = 106
$ws1.Range("F12").Value = 97
$ws1.Range("F13").Value = 781
$ws1.Range("F14").Value = 220
$ws1.Range("F18").Value = 340
$ws1.Range("F19").Value = 200
$ws1.Range("F20").Value = 690
$ws1.Range("F21").Value = 63
$ws1.Range("F23").Value = 179
$ws1.Range("F24").Value = 45
$ws1.Range("F25").Value = 898
$ws1.Range("F27").Value = 183
$ws1.Range("F28").Value = 56
$ws1.Range("F29").Value = 298
$ws1.Range("F30").Value = 15

# Sheet 2: 演出
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F11").Value = 128

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F3").Value = 1171
$ws4.Range("F4").Value = 1084
$ws4.Range("F5").Value = 1887
$ws4.Range("F6").Value = 592
$ws4.Range("F7").Value = 1235
$ws4.Range("F12").Value = 325
$ws4.Range("F13").Value = 106
$ws4.Range("F14").Value = 97
$ws4.Range("F15").Value = 781
$ws4.Range("F16").Value = 220
$ws4.Range("F23").Value = 340
$ws4.Range("F27").Value = 200
$ws4.Range("F28").Value = 690
$ws4.Range("F29").Value = 63
$ws4.Range("F31").Value = 179
$ws4.Range("F32").Value = 45
$ws4.Range("F33").Value = 898
$ws4.Range("F37").Value = 183
$ws4.Range("F38").Value = 56
$ws4.Range("F39").Value = 298
$ws4.Range("F41").Value = 128
$ws4.Range("F42").Value = 128
$ws4.Range("F43").Value = 15
